$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.10543283814906
$ws.Range("C2").Value = 10.91693200572579
$ws.Range("D2").Value = 4.025549227329006
$ws.Range("F2").Value = 17.68455412667648
$ws.Range("G2").Value = 3.59091937022017
$ws.Range("I2").Value = 15.46549118450062
$ws.Range("O2").Value = 15.96531596457525
$ws.Range("B3").Value = 14.25381441812363
$ws.Range("C3").Value = 10.29545017817815
$ws.Range("D3").Value = 3.916733090458262
$ws.Range("F3").Value = 17.80338077094061
$ws.Range("G3").Value = 3.592737965206635
$ws.Range("I3").Value = 15.63547274004895
$ws.Range("O3").Value = 16.11051273167243
$ws.Range("B4").Value = 13.70282031065429
$ws.Range("C4").Value = 9.892631557122902
$ws.Range("D4").Value = 3.847973633172853
$ws.Range("F4").Value = 17.88481039133162
$ws.Range("G4").Value = 3.593912627017501
$ws.Range("I4").Value = 15.74501813683899
$ws.Range("O4").Value = 16.20577188881545
$ws.Range("B5").Value = 13.47135605031164
$ws.Range("C5").Value = 9.723221134638369
$ws.Range("D5").Value = 3.81949373327941
$ws.Range("F5").Value = 17.92010582801657
$ws.Range("G5").Value = 3.594405949919424
$ws.Range("I5").Value = 15.79096423570918
$ws.Range("O5").Value = 16.24612103784295
$ws.Range("B6").Value = 13.43250761378878
$ws.Range("C6").Value = 9.694775739174849
$ws.Range("D6").Value = 3.814737796514466
$ws.Range("F6").Value = 17.92609366213862
$ws.Range("G6").Value = 3.594488751268458
$ws.Range("I6").Value = 15.79867251182958
$ws.Range("O6").Value = 16.25291324873805
$ws.Range("B7").Value = 13.69972655948722
$ws.Range("C7").Value = 9.890368002618809
$ws.Range("D7").Value = 3.847591364739379
$ws.Range("F7").Value = 17.88527786909745
$ws.Range("G7").Value = 3.593919220812297
$ws.Range("I7").Value = 15.74563249106744
$ws.Range("O7").Value = 16.20630986278315
$ws.Range("B8").Value = 14.81769866291513
$ws.Range("C8").Value = 10.70709240467504
$ws.Range("D8").Value = 3.988450287442977
$ws.Range("F8").Value = 17.72375723783713
$ws.Range("G8").Value = 3.591534403182032
$ws.Range("I8").Value = 15.52302889791615
$ws.Range("O8").Value = 16.01410873952816
$ws.Range("B9").Value = 16.78252483080188
$ws.Range("C9").Value = 12.13764797638032
$ws.Range("D9").Value = 4.248052387865902
$ws.Range("F9").Value = 17.47498602283876
$ws.Range("G9").Value = 3.587316229632771
$ws.Range("I9").Value = 15.12739336362547
$ws.Range("O9").Value = 15.68592206340534
$ws.Range("B10").Value = 18.08280603253113
$ws.Range("C10").Value = 13.08185620631934
$ws.Range("D10").Value = 4.427221484590204
$ws.Range("F10").Value = 17.33467288400696
$ws.Range("G10").Value = 3.584493733226143
$ws.Range("I10").Value = 14.86139244133518
$ws.Range("O10").Value = 15.47486766996478
$ws.Range("B11").Value = 18.64260802373132
$ws.Range("C11").Value = 13.48787908287064
$ws.Range("D11").Value = 4.505960779322924
$ws.Range("F11").Value = 17.28026951940991
$ws.Range("G11").Value = 3.583269151809175
$ws.Range("I11").Value = 14.7456847632316
$ws.Range("O11").Value = 15.38546514279299
$ws.Range("B12").Value = 18.84999635572084
$ws.Range("C12").Value = 13.63823211090986
$ws.Range("D12").Value = 4.535360681472108
$ws.Range("F12").Value = 17.26104016619515
$ws.Range("G12").Value = 3.582813928017558
$ws.Range("I12").Value = 14.70262698424277
$ws.Range("O12").Value = 15.35256831042284
$ws.Range("B13").Value = 18.80553649236036
$ws.Range("C13").Value = 13.60600231295288
$ws.Range("D13").Value = 4.529047717295553
$ws.Range("F13").Value = 17.26512027344328
$ws.Range("G13").Value = 3.582911591278184
$ws.Range("I13").Value = 14.71186658230693
$ws.Range("O13").Value = 15.35961049790801
$ws.Range("B14").Value = 18.65976228773348
$ws.Range("C14").Value = 13.50031692789402
$ws.Range("D14").Value = 4.508387986991902
$ws.Range("F14").Value = 17.27865991895068
$ws.Range("G14").Value = 3.583231530183558
$ws.Range("I14").Value = 14.74212719854345
$ws.Range("O14").Value = 15.38273945551171
$ws.Range("B15").Value = 18.56987183007436
$ws.Range("C15").Value = 13.43513845388434
$ws.Range("D15").Value = 4.495678466217673
$ws.Range("F15").Value = 17.28713250144471
$ws.Range("G15").Value = 3.583428607557094
$ws.Range("I15").Value = 14.76076133547161
$ws.Range("O15").Value = 15.39703160333164
$ws.Range("B16").Value = 18.04558025180255
$ws.Range("C16").Value = 13.05484692556494
$ws.Range("D16").Value = 4.422018343854603
$ws.Range("F16").Value = 17.33841937739593
$ws.Range("G16").Value = 3.584574953989031
$ws.Range("I16").Value = 14.86906050375885
$ws.Range("O16").Value = 15.48084395401528
$ws.Range("B17").Value = 17.71579429489122
$ws.Range("C17").Value = 12.81551502654688
$ws.Range("D17").Value = 4.376107574826038
$ws.Range("F17").Value = 17.37230840836671
$ws.Range("G17").Value = 3.585293381343683
$ws.Range("I17").Value = 14.93685280234017
$ws.Range("O17").Value = 15.53395787458071
$ws.Range("B18").Value = 17.52312765741285
$ws.Range("C18").Value = 12.67564622229494
$ws.Range("D18").Value = 4.349441625939582
$ws.Range("F18").Value = 17.39268672898159
$ws.Range("G18").Value = 3.585712193916136
$ws.Range("I18").Value = 14.97634397409608
$ws.Range("O18").Value = 15.5651290074761
$ws.Range("B19").Value = 17.45738312684662
$ws.Range("C19").Value = 12.62790993142101
$ws.Range("D19").Value = 4.340369079612126
$ws.Range("F19").Value = 17.39973817506346
$ws.Range("G19").Value = 3.585854958373241
$ws.Range("I19").Value = 14.9898007906843
$ws.Range("O19").Value = 15.57578951056198
$ws.Range("B20").Value = 17.75120955484799
$ws.Range("C20").Value = 12.84122135072608
$ws.Range("D20").Value = 4.381021835670222
$ws.Range("F20").Value = 17.3686090175752
$ws.Range("G20").Value = 3.585216325056208
$ws.Range("I20").Value = 14.92958459478561
$ws.Range("O20").Value = 15.52823943735401
$ws.Range("B21").Value = 18.70270470377172
$ws.Range("C21").Value = 13.53145166018312
$ws.Range("D21").Value = 4.514467710771606
$ws.Range("F21").Value = 17.27464563169832
$ws.Range("G21").Value = 3.583137326047663
$ws.Range("I21").Value = 14.73321837274998
$ws.Range("O21").Value = 15.37591985705403
$ws.Range("B22").Value = 19.29775903970291
$ws.Range("C22").Value = 13.96274006439345
$ws.Range("D22").Value = 4.599244339004692
$ws.Range("F22").Value = 17.22123951356166
$ws.Range("G22").Value = 3.581828097594416
$ws.Range("I22").Value = 14.60929947505194
$ws.Range("O22").Value = 15.28195780842555
$ws.Range("B23").Value = 18.98262969364763
$ws.Range("C23").Value = 13.73437146063489
$ws.Range("D23").Value = 4.554226404562328
$ws.Range("F23").Value = 17.24900565738522
$ws.Range("G23").Value = 3.58252234008862
$ws.Range("I23").Value = 14.67503426582817
$ws.Range("O23").Value = 15.33159318410413
$ws.Range("B24").Value = 17.73520786994236
$ws.Range("C24").Value = 12.82960660217913
$ws.Range("D24").Value = 4.378800941703216
$ws.Range("F24").Value = 17.37027872535326
$ws.Range("G24").Value = 3.58525114419854
$ws.Range("I24").Value = 14.93286894232612
$ws.Range("O24").Value = 15.53082276468477
$ws.Range("B25").Value = 16.27588966320838
$ws.Range("C25").Value = 11.76926332468754
$ws.Range("D25").Value = 4.179758405237664
$ws.Range("F25").Value = 17.53490124072817
$ws.Range("G25").Value = 3.588408576901134
$ws.Range("I25").Value = 15.23007220239142
$ws.Range("O25").Value = 15.76945294430498
